$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 9.207199999999997
$ws.Range("B6").Value = 6.571999999999999
$ws.Range("B7").Value = 5.554199999999998
$ws.Range("C7").Value = -14.0701
$ws.Range("B8").Value = 6.397400000000001
$ws.Range("C11").Value = -12.6207
$ws.Range("C12").Value = -11.38099999999999
$ws.Range("E12").Value = 17.42840000000002
$ws.Range("E13").Value = 16.8296
$ws.Range("E14").Value = 16.8918
$ws.Range("C15").Value = -14.24939999999999
$ws.Range("B16").Value = 5.221399999999999
$ws.Range("E16").Value = 16.0982
$ws.Range("E19").Value = 16.47649999999999
$ws.Range("B20").Value = 9.841899999999988
$ws.Range("C20").Value = -12.5
$ws.Range("E20").Value = 15.92489999999999
$ws.Range("B21").Value = 8.884299999999989
$ws.Range("C21").Value = -12.0764
$ws.Range("C22").Value = -12.0437
$ws.Range("E22").Value = 16.91760000000002
$ws.Range("C23").Value = -12.1433
$ws.Range("B28").Value = 5.835600000000002
$ws.Range("B29").Value = 4.774099999999998
$ws.Range("C29").Value = -10.29690000000001
$ws.Range("B30").Value = 4.784599999999998
$ws.Range("B32").Value = 7.541299999999995
$ws.Range("C34").Value = -11.63830000000001
$ws.Range("E36").Value = 15.7253
$ws.Range("B40").Value = 9.271799999999995
$ws.Range("C42").Value = -12.28930000000001
$ws.Range("C43").Value = -12.8926
$ws.Range("E43").Value = 17.20150000000001
$ws.Range("C44").Value = -14.27329999999999
$ws.Range("C45").Value = -13.75739999999999
$ws.Range("B46").Value = 5.761599999999999
$ws.Range("C46").Value = -13.78399999999999
$ws.Range("E46").Value = 16.2842
$ws.Range("C50").Value = -13.98789999999999
$ws.Range("E50").Value = 16.73889999999999
$ws.Range("B51").Value = 6.284000000000002
$ws.Range("C51").Value = -13.12390000000001
$ws.Range("B52").Value = 5.880100000000001
$ws.Range("B57").Value = 5.748899999999999
$ws.Range("C57").Value = -13.68649999999999
$ws.Range("B59").Value = 6.693400000000002
$ws.Range("B62").Value = 6.072399999999996
$ws.Range("C65").Value = -13.1659
$ws.Range("B66").Value = 5.736599999999998
$ws.Range("C66").Value = -11.738
$ws.Range("C67").Value = -11.1672
$ws.Range("B73").Value = 8.358900000000002
$ws.Range("B74").Value = 9.30139999999999
$ws.Range("E76").Value = 16.24419999999998
$ws.Range("B77").Value = 8.946300000000001
$ws.Range("C79").Value = -11.27330000000001
$ws.Range("C84").Value = -13.46569999999999
$ws.Range("C87").Value = -13.7684
$ws.Range("B92").Value = 5.533299999999994
$ws.Range("C92").Value = -11.311
$ws.Range("E95").Value = 18.26320000000002
$ws.Range("C97").Value = -11.2119
$ws.Range("E97").Value = 16.5574
$ws.Range("E99").Value = 16.30209999999998
$ws.Range("B100").Value = 6.3714
